# Generate Report for Handoff
#
# 2caa79a4-97ad-499e-a178-6fe2978ed58a.md finished its handoff ("Ready for
# handoff") while 790442b1-f051-4ee8-8be6-9e7016305ede.md is still
# "In Translation" -- the status report's row ordering / timestamps are
# regenerated, which swaps the two files' rows 9 & 10 (on every language
# sheet) and refreshes a couple of "Latest Handoff/Handback" timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Refresh handback timestamp for 5e33b3fc row
$ws.Range("D6").Value = "2016-23-18 22:23:49"

# Row 9 now reports 2caa79a4 (ready for handoff)
$ws.Range("A9").Hyperlinks.Item(1).TextToDisplay = "2caa79a4-97ad-499e-a178-6fe2978ed58a.md"
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "2016-23-18 22:23:49"

# Row 10 now reports 790442b1 (still in translation, but keeps "Ready for
# handoff" cached status like the generator emitted)
$ws.Range("A10").Hyperlinks.Item(1).TextToDisplay = "790442b1-f051-4ee8-8be6-9e7016305ede.md"
$ws.Range("D10").Value = "2016-23-18 22:23:49"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("E6").Value = "2016-03-18 22:23:43"

$ws.Range("A9").Hyperlinks.Item(1).TextToDisplay = "2caa79a4-97ad-499e-a178-6fe2978ed58a.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Hyperlinks.Item(1).TextToDisplay = "2caa79a4-97ad-499e-a178-6fe2978ed58a.d7ab2f752541d7b8e5dcf93bd932c789bc177340.zh-cn.xlf"
$ws.Range("E9").Value = "2016-03-18 22:23:43"

$ws.Range("A10").Hyperlinks.Item(1).TextToDisplay = "790442b1-f051-4ee8-8be6-9e7016305ede.md"
$ws.Range("D10").Hyperlinks.Item(1).TextToDisplay = "790442b1-f051-4ee8-8be6-9e7016305ede.7344ed57d998fcac44a689258e4a5d75b79b3ab8.zh-cn.xlf"
$ws.Range("E10").Value = "2016-03-18 22:23:43"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("E6").Value = "2016-03-18 22:23:49"

$ws.Range("A9").Hyperlinks.Item(1).TextToDisplay = "2caa79a4-97ad-499e-a178-6fe2978ed58a.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Hyperlinks.Item(1).TextToDisplay = "2caa79a4-97ad-499e-a178-6fe2978ed58a.d7ab2f752541d7b8e5dcf93bd932c789bc177340.de-de.xlf"
$ws.Range("E9").Value = "2016-03-18 22:23:49"

$ws.Range("A10").Hyperlinks.Item(1).TextToDisplay = "790442b1-f051-4ee8-8be6-9e7016305ede.md"
$ws.Range("D10").Hyperlinks.Item(1).TextToDisplay = "790442b1-f051-4ee8-8be6-9e7016305ede.7344ed57d998fcac44a689258e4a5d75b79b3ab8.de-de.xlf"
$ws.Range("E10").Value = "2016-03-18 22:23:49"
